# Daily attendance processing - 2025-10-16 19:43:27
# Rotate the comma-separated "Recorded By" entries in column G so the
# last entry moves to the front (applies to every data row on the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",") | ForEach-Object { $_.Trim() }
        $count = $parts.Count

        if ($count -gt 1) {
            $rotated = @($parts[$count - 1]) + $parts[0..($count - 2)]
            $cell.Value2 = [string]::Join(", ", $rotated)
        }
    }
}
